$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'30.158.68"
$ws.Range('E2').Value = '  -1.46%  '
$ws.Range('D3').Value = "'1.853.80"
$ws.Range('E3').Value = '  -2.11%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'235.53"
$ws.Range('E5').Value = '  -1.28%  '
$ws.Range('D6').Value = "'1.002"
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').Value = "'0.4765"
$ws.Range('E7').Value = '  -2.51%  '
$ws.Range('D8').Value = "'0.2810"
$ws.Range('E8').Value = '  -4.17%  '
$ws.Range('D9').Value = "'0.06501"
$ws.Range('E9').Value = '  -2.96%  '
$ws.Range('D10').Value = "'1.846.54"
$ws.Range('E10').Value = '  -2.36%  '
$ws.Range('D11').Value = "'0.07344"
$ws.Range('E11').Value = '  +0.06%  '
$ws.Range('D12').Value = "'16.31"
$ws.Range('E12').Value = '  -4.84%  '
$ws.Range('D13').Value = "'5.133"
$ws.Range('E13').Value = '  -0.34%  '
$ws.Range('D14').Value = "'87.18"
$ws.Range('E14').Value = '  -1.05%  '
$ws.Range('D15').Value = "'0.6435"
$ws.Range('E15').Value = '  -3.57%  '
$ws.Range('D16').Value = "'30.103.36"
$ws.Range('E16').Value = '  -1.43%  '
$ws.Range('D17').Value = "'13.25"
$ws.Range('E17').Value = '  -1.35%  '
$ws.Range('D18').Value = "'1.002"
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').Value = "'0.000007625"
$ws.Range('E19').Value = '  -2.84%  '
$ws.Range('D20').Value = "'2.123.77"
$ws.Range('E20').Value = '  -0.53%  '
$ws.Range('D21').Value = "'1.002"
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('D22').Value = "'5.252"
$ws.Range('E22').Value = '  -0.86%  '
$ws.Range('D23').Value = "'218.46"
$ws.Range('E23').Value = '  +15.27%  '
$ws.Range('D24').Value = "'6.102"
$ws.Range('E24').Value = '  -1.18%  '
$ws.Range('D25').Value = "'9.276"
$ws.Range('E25').Value = '  -2.02%  '
$ws.Range('D26').Value = "'165.13"
$ws.Range('E26').Value = '  +1.90%  '
$ws.Range('D27').Value = "'18.53"
$ws.Range('E27').Value = '  +1.16%  '
$ws.Range('D28').Value = "'1.906"
$ws.Range('E28').Value = '  -1.27%  '
$ws.Range('D29').Value = "'1.427"
$ws.Range('E29').Value = '  -3.26%  '
$ws.Range('D30').Value = "'4.250"
$ws.Range('E30').Value = '  -2.90%  '
$ws.Range('D31').Value = "'0.09140"
$ws.Range('E31').Value = '  -0.13%  '
$ws.Range('D32').Value = "'3.969"
$ws.Range('E32').Value = '  -3.35%  '
$ws.Range('D33').Value = "'0.05015"
$ws.Range('E33').Value = '  -3.90%  '
$ws.Range('D34').Value = "'0.7411"
$ws.Range('E34').Value = '  +0.41%  '
$ws.Range('D35').Value = "'1.138"
$ws.Range('E35').Value = '  +3.42%  '
$ws.Range('D36').Value = "'2.692"
$ws.Range('E36').Value = '  -0.87%  '
$ws.Range('D37').Value = "'0.01820"
$ws.Range('E37').Value = '  -0.58%  '
$ws.Range('D38').Value = "'2.615"
$ws.Range('E38').Value = '  -2.60%  '
$ws.Range('D39').Value = "'0.9004"
$ws.Range('E39').Value = '  -2.22%  '
$ws.Range('D40').Value = "'2.037"
$ws.Range('E40').Value = '  -0.84%  '
$ws.Range('D41').Value = "'5.936"
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('D42').Value = "'106.49"
$ws.Range('E42').Value = '  +0.20%  '
$ws.Range('D43').Value = "'0.4243"
$ws.Range('E43').Value = '  -3.63%  '
$ws.Range('D44').Value = "'1.001"
$ws.Range('E44').Value = '  +0.73%  '
$ws.Range('D45').Value = "'7.419"
$ws.Range('E45').Value = '  -2.21%  '
$ws.Range('D46').Value = "'0.1309"
$ws.Range('E46').Value = '  -5.43%  '
$ws.Range('D47').Value = "'1.558"
$ws.Range('E47').Value = '  +9.54%  '
$ws.Range('D48').Value = "'63.96"
$ws.Range('E48').Value = '  -7.24%  '
$ws.Range('D49').Value = "'8.812"
$ws.Range('E49').Value = '  -1.85%  '
$ws.Range('D50').Value = "'34.19"
$ws.Range('E50').Value = '  -2.31%  '
$ws.Range('D51').Value = "'0.05687"
$ws.Range('E51').Value = '  -2.37%  '
